$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add new column headers D1=3, E1=4 (same style as C1)
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("D1").Style = $ws.Range("C1").Style

# Update existing column C values (rows 2-13)
$ws.Range("C2").Value = -5.162938253646961
$ws.Range("C3").Value = -1.230872123270181
$ws.Range("C4").Value = -0.08073787866549097
$ws.Range("C5").Value = -0.4375737637739642
$ws.Range("C6").Value = 0.01336035328125892
$ws.Range("C7").Value = 0.09957005506681453
$ws.Range("C8").Value = 0.1255627463891133
$ws.Range("C9").Value = 0.03065774420266123
$ws.Range("C10").Value = 0.02915937565228036
$ws.Range("C11").Value = 0.003142978813154578
$ws.Range("C12").Value = 0.0416495879320569
$ws.Range("C13").Value = -0.001427841317002163

# Fill new column D (rows 2-13)
$ws.Range("D2").Value = -4.814467328962087
$ws.Range("D3").Value = -1.212117344132834
$ws.Range("D4").Value = 0.01491559983289693
$ws.Range("D5").Value = -0.1734273993507124
$ws.Range("D6").Value = -0.04886876058402497
$ws.Range("D7").Value = 0.06232481319395875
$ws.Range("D8").Value = 0.03018733153560915
$ws.Range("D9").Value = 0.0255049138530617
$ws.Range("D10").Value = 0.02538142198961245
$ws.Range("D11").Value = 0.008526450440452422
$ws.Range("D12").Value = 0.02740891658032532
$ws.Range("D13").Value = 0.001083974367374332

# Fill new column E (rows 2-13)
$ws.Range("E2").Value = -4.434346894216289
$ws.Range("E3").Value = -1.181443779698794
$ws.Range("E4").Value = 0.08984950567291528
$ws.Range("E5").Value = 0.03793664658911642
$ws.Range("E6").Value = -0.1037897773736463
$ws.Range("E7").Value = 0.02488347623765172
$ws.Range("E8").Value = -0.0607292752145601
$ws.Range("E9").Value = 0.02259250009160615
$ws.Range("E10").Value = 0.01743950579382542
$ws.Range("E11").Value = 0.01036837390730543
$ws.Range("E12").Value = 0.01773775477520047
$ws.Range("E13").Value = 0.002276948278055886
